# agrego preguntas de siniestro motor para aprobaciones previas a pago
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$claimCenterUrl = "https://preproducciongestion.segurossura.com.ar/cc/ClaimCenter.do"

# Insert a new row above current row 2; rows 2-5 shift down to 3-6.
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with the "motor" claim entry (mgentilini / 0420194406715).
$ws.Range("A2").Value = "ssurgwsoadev4-oci.opc.oracleoutsourcing.com"
$ws.Range("B2").Value = "https://ssurgwsoadev4-oci.opc.oracleoutsourcing.com/cc/ClaimCenter.do"
$ws.Range("E2").Value = "'0420194406715 "
$ws.Range("C2").Value = "mgentilini"
$ws.Range("D2").Value = "silverarrow"

# The row insert doesn't carry the hyperlinks on B4/B5 down to B5/B6 with it, so
# rebuild them on the now-shifted cells, preserving the original cell styling.
$linkStyle = $ws.Range("B5").Style
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B6"), $claimCenterUrl) | Out-Null
$ws.Hyperlinks.Add($ws.Range("B5"), $claimCenterUrl) | Out-Null
$ws.Range("B5").Style = $linkStyle
$ws.Range("B6").Style = $linkStyle

# Move the selection, matching the recorded UI state after the edit.
$ws.Range("I14").Select()
